$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 196.85715
$ws.Range("I9").Value = 196.85715
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 196.85715
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = -27.85714999999999
$ws.Range("N9").Value = ""

$ws.Range("H86").Value = 45848570
$ws.Range("I86").Value = 68771660
$ws.Range("J86").Value = 2369.25
$ws.Range("K86").Value = 68771660
$ws.Range("L86").Value = 2369.25
$ws.Range("M86").Value = -68770537

$ws.Range("H89").Value = 45848570
$ws.Range("I89").Value = 68771660
$ws.Range("J89").Value = 2369.25
$ws.Range("K89").Value = 343858300
$ws.Range("L89").Value = 11846.25
$ws.Range("M89").Value = -343852684

$ws.Range("H98").Value = 1387.8096
$ws.Range("I98").Value = 1097.4117
$ws.Range("J98").Value = 2622
$ws.Range("K98").Value = 1097.4117
$ws.Range("L98").Value = 2622
$ws.Range("M98").Value = 400.5882999999999
$ws.Range("N98").Value = -5618

$ws.Range("H122").Value = 1387.8096
$ws.Range("I122").Value = 1097.4117
$ws.Range("J122").Value = 2622
$ws.Range("K122").Value = 3292.2351
$ws.Range("L122").Value = 7866
$ws.Range("M122").Value = -842.2351000000003
$ws.Range("N122").Value = -12766

$ws.Range("H132").Value = 4224.522
$ws.Range("I132").Value = 4007.8572
$ws.Range("J132").Value = 6499.5
$ws.Range("K132").Value = 12023.5716
$ws.Range("L132").Value = 19498.5
$ws.Range("M132").Value = -9493.571599999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H21").Value = 40209.8
$ws.Range("I21").Value = 10338.333
$ws.Range("J21").Value = 85017
$ws.Range("K21").Value = 10338.333
$ws.Range("L21").Value = 85017
$ws.Range("M21").Value = -9964.333000000001
$ws.Range("N21").Value = -85765

$ws.Range("H32").Value = 1801992
$ws.Range("I32").Value = 2200320.8
$ws.Range("J32").Value = 21228.295
$ws.Range("K32").Value = 2200320.8
$ws.Range("L32").Value = 21228.295
$ws.Range("M32").Value = -2200033.8

$ws.Range("H132").Value = 1879765.2
$ws.Range("I132").Value = 2923.8965
$ws.Range("J132").Value = 6415465.5
$ws.Range("K132").Value = 8771.6895
$ws.Range("L132").Value = 19246396.5
$ws.Range("M132").Value = -6241.6895

$ws.Range("H133").Value = 47505.5
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 47505.5
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 47505.5
$ws.Range("N133").Value = -52565.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2236.6667
$ws.Range("I99").Value = 2236.6667
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 2236.6667
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -738.6667000000002
$ws.Range("N99").Value = ""

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H35").Value = 501750
$ws.Range("I35").Value = 501750
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 501750
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -501456

$ws.Range("H41").Value = 12392.6
$ws.Range("I41").Value = 5900
$ws.Range("J41").Value = 14015.75
$ws.Range("K41").Value = 5900
$ws.Range("L41").Value = 14015.75
$ws.Range("M41").Value = -5472
$ws.Range("N41").Value = -14871.75

$ws.Range("H50").Value = 13665.833
$ws.Range("I50").Value = 10000
$ws.Range("J50").Value = 14399
$ws.Range("K50").Value = 10000
$ws.Range("L50").Value = 14399
$ws.Range("M50").Value = -9375
$ws.Range("N50").Value = -15649

$ws.Range("H51").Value = 18042
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 18042
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 18042
$ws.Range("N51").Value = -19514

$ws.Range("H59").Value = 14948.857
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 14948.857
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 14948.857
$ws.Range("N59").Value = -17238.857

$ws.Range("H60").Value = 15134.083
$ws.Range("I60").Value = 7733.3335
$ws.Range("J60").Value = 17601
$ws.Range("K60").Value = 7733.3335
$ws.Range("L60").Value = 17601
$ws.Range("M60").Value = -7222.3335
$ws.Range("N60").Value = -18623

$ws.Range("H61").Value = 18042
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 18042
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 18042
$ws.Range("N61").Value = -18738

$ws.Range("H68").Value = 22090.363
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 22090.363
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 22090.363
$ws.Range("N68").Value = -23588.363

$ws.Range("H71").Value = 22090.363
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 22090.363
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 66271.08900000001
$ws.Range("N71").Value = -73759.08900000001

$ws.Range("H74").Value = 18833
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 18833
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 18833
$ws.Range("N74").Value = -20581

$ws.Range("H77").Value = 18833
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 18833
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 56499
$ws.Range("N77").Value = -65235

$ws.Range("H94").Value = 76924360
$ws.Range("I94").Value = 200000720
$ws.Range("J94").Value = 1628.5
$ws.Range("K94").Value = 200000720
$ws.Range("L94").Value = 1628.5
$ws.Range("M94").Value = -200000269
$ws.Range("N94").Value = -2530.5

$ws.Range("H122").Value = 1547.6207
$ws.Range("I122").Value = 1267.2
$ws.Range("J122").Value = 1695.2106
$ws.Range("K122").Value = 3801.6
$ws.Range("L122").Value = 5085.6318
$ws.Range("M122").Value = -1351.6
$ws.Range("N122").Value = -9985.631799999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 271
$ws.Range("I11").Value = 125
$ws.Range("J11").Value = 297.54544
$ws.Range("K11").Value = 375
$ws.Range("L11").Value = 892.63632
$ws.Range("M11").Value = -235
$ws.Range("N11").Value = -1172.63632

$ws.Range("H98").Value = 317066.66
$ws.Range("I98").Value = 700
$ws.Range("J98").Value = 475250
$ws.Range("K98").Value = 2100
$ws.Range("L98").Value = 1425750
$ws.Range("M98").Value = -602
$ws.Range("N98").Value = -1428746

$ws.Range("H139").Value = 3982.238
$ws.Range("I139").Value = 1415.6364
$ws.Range("J139").Value = 6805.5
$ws.Range("K139").Value = 4246.9092
$ws.Range("L139").Value = 20416.5
$ws.Range("M139").Value = 893.0907999999999
$ws.Range("N139").Value = -30696.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").Value = ""

$ws.Range("H102").Value = 1499.5714
$ws.Range("I102").Value = 1499.5714
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 1499.5714
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = 122.4286

$ws.Range("H113").Value = 80543.78999999999
$ws.Range("I113").Value = 93617.75
$ws.Range("J113").Value = 2100
$ws.Range("K113").Value = 93617.75
$ws.Range("L113").Value = 2100
$ws.Range("M113").Value = -91447.75
$ws.Range("N113").Value = -6440

$ws.Range("H132").Value = 2736.75
$ws.Range("I132").Value = 2361.2778
$ws.Range("J132").Value = 3412.6
$ws.Range("K132").Value = 7083.8334
$ws.Range("L132").Value = 10237.8
$ws.Range("M132").Value = -4553.8334
$ws.Range("N132").Value = -15297.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8099.1113
$ws.Range("I7").Value = 8270.286
$ws.Range("J7").Value = 7500
$ws.Range("K7").Value = 8270.286
$ws.Range("L7").Value = 7500
$ws.Range("M7").Value = -8158.286
$ws.Range("N7").Value = -7724

$ws.Range("H40").Value = 1776.4667
$ws.Range("I40").Value = 1617.6428
$ws.Range("J40").Value = 4000
$ws.Range("K40").Value = 1617.6428
$ws.Range("L40").Value = 4000
$ws.Range("M40").Value = -1481.6428
$ws.Range("N40").Value = -4272

$ws.Range("H41").Value = 2033
$ws.Range("I41").Value = 2033
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 2033
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = -1595

$ws.Range("H88").Value = 38500

$ws.Range("H91").Value = 38500

$ws.Range("H98").Value = 92500
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 92500
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 92500
$ws.Range("N98").Value = -98490

$ws.Range("H126").Value = 8099.1113
$ws.Range("I126").Value = 8270.286
$ws.Range("J126").Value = 7500
$ws.Range("K126").Value = 24810.858
$ws.Range("L126").Value = 22500
$ws.Range("M126").Value = -22340.858
$ws.Range("N126").Value = -27440

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H101").Value = 88602
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 88602
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 88602
$ws.Range("N101").Value = -95092
